$d = $word.ActiveDocument

# Locate the paragraph that contains the "Expression ... is invalid" error
# message (the second paragraph of the document) and compute the insertion
# point just before its trailing paragraph mark.
$p = $d.Paragraphs(2)
$pos = $p.Range.End - 1

function Insert-PlainRun($position, [string]$text) {
    $r = $d.Range($position, $position)
    $r.InsertAfter($text)
    return $r.End
}

function Insert-ErrorRun($position, [string]$text) {
    $r = $d.Range($position, $position)
    $r.InsertAfter($text)
    $r.Font.Color = 255              # wdColorRed
    $r.Font.Size = 16                # 16pt -> w:sz 32 (half-points)
    $r.Font.HighlightColorIndex = 16 # wdGray25 -> w:highlight lightGray
    return $r.End
}

$pos = Insert-PlainRun $pos "    "
$pos = Insert-ErrorRun $pos "<---"
$pos = Insert-ErrorRun $pos "Couldn't find the 'self' variable"
$pos = Insert-PlainRun $pos "    "
$pos = Insert-ErrorRun $pos "<---"
$pos = Insert-ErrorRun $pos "missing feature access or service call"
